$wb = $excel.ActiveWorkbook

# --- Add "Romania" sheet (copied from Spain, placed at the end) ---
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$romania = $wb.Worksheets.Item($wb.Worksheets.Count)
$romania.Name = "Romania"
$romania.Range("B2").Value = "Romania Market"
$romania.Range("B4").Value = "NGC-4307/T3532/T3545"
$romania.Range("B8").Select()

# --- Add "Slovakia" sheet (copied from Spain, placed at the end) ---
$spain2 = $wb.Worksheets.Item("Spain")
$spain2.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"
$slovakia.Range("B4").Value = "NGC-4306/T3565/T3577"
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B3").Select()

# Slovakia is the last sheet and becomes the active/selected tab.
$slovakia.Activate()
